# Generate ML input with median values and merge with site points
#
# The "workshop_ml_input" column (C) on the "workshop_feature_selection"
# sheet is a binary 0/1 flag. This run flips every remaining 0 to 1 across
# the data rows (C2:C83), marking every feature as included in the ML input.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workshop_feature_selection")

for ($r = 2; $r -le 83; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 0) {
        $cell.Value = 1
    }
}
